$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“개발자를 위한 필수 수학”이 출간되었습니다!"
$ws.Range("E12").Value = "https://tensorflow.blog/2024/05/28/%ea%b0%9c%eb%b0%9c%ec%9e%90%eb%a5%bc-%ec%9c%84%ed%95%9c-%ed%95%84%ec%88%98-%ec%88%98%ed%95%99%ec%9d%b4-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D24").Value = "Daily arXiv: 2024-05-27"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/223460598728?fromRss=true&trackingCode=rss"

$ws.Range("D51").Value = "맥북, 맥미니 화면 잠금 단축키"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%EB%A7%A5%EB%B6%81-%EB%A7%A5%EB%AF%B8%EB%8B%88-%ED%99%94%EB%A9%B4-%EC%9E%A0%EA%B8%88-%EB%8B%A8%EC%B6%95%ED%82%A4"
